$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 38.93002066666667
$ws.Range("H2").Value = 116.790062
$ws.Range("I2").Value = 0.7610372167397395
$ws.Range("J2").Value = 0.7610372167397393
$ws.Range("M2").Value = 71.541692
$ws.Range("N2").Value = 214.625076
$ws.Range("O2").Value = 0.6133071420247926
$ws.Range("P2").Value = 0.6133071420247926
$ws.Range("Q2").Value = 2785.119548088302
$ws.Range("R2").Value = 25066.07593279471
$ws.Range("S2").Value = 0.4667495603731523
$ws.Range("T2").Value = 0.4667495603731522

$ws.Range("G3").Value = 38.93002066666667
$ws.Range("H3").Value = 116.790062
$ws.Range("I3").Value = 0.7610372167397395
$ws.Range("J3").Value = 0.7610372167397393
$ws.Range("O3").Value = 0.08457024278578675
$ws.Range("P3").Value = 0.08457024278578675
$ws.Range("Q3").Value = 384.0461332174521
$ws.Range("R3").Value = 3456.415198957068
$ws.Range("S3").Value = 0.06436110218869917
$ws.Range("T3").Value = 0.06436110218869917

$ws.Range("G4").Value = 38.93002066666667
$ws.Range("H4").Value = 116.790062
$ws.Range("I4").Value = 0.7610372167397395
$ws.Range("J4").Value = 0.7610372167397393
$ws.Range("M4").Value = 35.05835333333334
$ws.Range("N4").Value = 105.17506
$ws.Range("O4").Value = 0.3005455684073286
$ws.Range("P4").Value = 0.3005455684073286
$ws.Range("Q4").Value = 1364.822419805969
$ws.Range("R4").Value = 12283.40177825372
$ws.Range("S4").Value = 0.2287263628841763
$ws.Range("T4").Value = 0.2287263628841763

$ws.Range("G5").Value = 38.93002066666667
$ws.Range("H5").Value = 116.790062
$ws.Range("I5").Value = 0.7610372167397395
$ws.Range("J5").Value = 0.7610372167397393
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.183961
$ws.Range("N5").Value = 0.551883
$ws.Range("O5").Value = 0.001577046782092083
$ws.Range("P5").Value = 0.001577046782092083
$ws.Range("Q5").Value = 7.161605531860668
$ws.Range("R5").Value = 64.45444978674601
$ws.Range("S5").Value = 0.001200191293711721
$ws.Range("T5").Value = 0.001200191293711721

$ws.Range("I6").Value = 0.1914142145281647
$ws.Range("J6").Value = 0.1914142145281647
$ws.Range("M6").Value = 71.541692
$ws.Range("N6").Value = 214.625076
$ws.Range("O6").Value = 0.6133071420247926
$ws.Range("P6").Value = 0.6133071420247926
$ws.Range("Q6").Value = 700.5064390256666
$ws.Range("R6").Value = 6304.557951230999
$ws.Range("S6").Value = 0.1173957048551893
$ws.Range("T6").Value = 0.1173957048551892

$ws.Range("I7").Value = 0.1914142145281647
$ws.Range("J7").Value = 0.1914142145281647
$ws.Range("O7").Value = 0.08457024278578675
$ws.Range("P7").Value = 0.08457024278578675
$ws.Range("S7").Value = 0.01618794659529756
$ws.Range("T7").Value = 0.01618794659529756

$ws.Range("I8").Value = 0.1914142145281647
$ws.Range("J8").Value = 0.1914142145281647
$ws.Range("M8").Value = 35.05835333333334
$ws.Range("N8").Value = 105.17506
$ws.Range("O8").Value = 0.3005455684073286
$ws.Range("P8").Value = 0.3005455684073286
$ws.Range("Q8").Value = 343.2767881927778
$ws.Range("R8").Value = 3089.491093735
$ws.Range("S8").Value = 0.05752869390660961
$ws.Range("T8").Value = 0.0575286939066096

$ws.Range("I9").Value = 0.1914142145281647
$ws.Range("J9").Value = 0.1914142145281647
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.183961
$ws.Range("N9").Value = 0.551883
$ws.Range("O9").Value = 0.001577046782092083
$ws.Range("P9").Value = 0.001577046782092083
$ws.Range("Q9").Value = 1.801269461583334
$ws.Range("R9").Value = 16.21142515425
$ws.Range("S9").Value = 0.0003018691710683258
$ws.Range("T9").Value = 0.0003018691710683257

$ws.Range("G10").Value = 1.794146
$ws.Range("H10").Value = 5.382438
$ws.Range("I10").Value = 0.03507349482179579
$ws.Range("J10").Value = 0.03507349482179579
$ws.Range("M10").Value = 71.541692
$ws.Range("N10").Value = 214.625076
$ws.Range("O10").Value = 0.6133071420247926
$ws.Range("P10").Value = 0.6133071420247926
$ws.Range("Q10").Value = 128.356240535032
$ws.Range("R10").Value = 1155.206164815288
$ws.Range("S10").Value = 0.02151082486997694
$ws.Range("T10").Value = 0.02151082486997694

$ws.Range("G11").Value = 1.794146
$ws.Range("H11").Value = 5.382438
$ws.Range("I11").Value = 0.03507349482179579
$ws.Range("J11").Value = 0.03507349482179579
$ws.Range("O11").Value = 0.08457024278578675
$ws.Range("P11").Value = 0.08457024278578675
$ws.Range("Q11").Value = 17.699318467548
$ws.Range("R11").Value = 159.293866207932
$ws.Range("S11").Value = 0.002966173972425305
$ws.Range("T11").Value = 0.002966173972425304

$ws.Range("G12").Value = 1.794146
$ws.Range("H12").Value = 5.382438
$ws.Range("I12").Value = 0.03507349482179579
$ws.Range("J12").Value = 0.03507349482179579
$ws.Range("M12").Value = 35.05835333333334
$ws.Range("N12").Value = 105.17506
$ws.Range("O12").Value = 0.3005455684073286
$ws.Range("P12").Value = 0.3005455684073286
$ws.Range("Q12").Value = 62.89980439958668
$ws.Range("R12").Value = 566.0982395962801
$ws.Range("S12").Value = 0.01054118343724811
$ws.Range("T12").Value = 0.01054118343724811

$ws.Range("G13").Value = 1.794146
$ws.Range("H13").Value = 5.382438
$ws.Range("I13").Value = 0.03507349482179579
$ws.Range("J13").Value = 0.03507349482179579
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.183961
$ws.Range("N13").Value = 0.551883
$ws.Range("O13").Value = 0.001577046782092083
$ws.Range("P13").Value = 0.001577046782092083
$ws.Range("Q13").Value = 0.330052892306
$ws.Range("R13").Value = 2.970476030754
$ws.Range("S13").Value = 0.00005531254214543638
$ws.Range("T13").Value = 0.00005531254214543637

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6381486666666666
$ws.Range("H14").Value = 1.914446
$ws.Range("I14").Value = 0.01247507391030006
$ws.Range("J14").Value = 0.01247507391030006
$ws.Range("M14").Value = 71.541692
$ws.Range("N14").Value = 214.625076
$ws.Range("O14").Value = 0.6133071420247926
$ws.Range("P14").Value = 0.6133071420247926
$ws.Range("Q14").Value = 45.65423536087732
$ws.Range("R14").Value = 410.8881182478959
$ws.Range("S14").Value = 0.007651051926474187
$ws.Range("T14").Value = 0.007651051926474187

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6381486666666666
$ws.Range("H15").Value = 1.914446
$ws.Range("I15").Value = 0.01247507391030006
$ws.Range("J15").Value = 0.01247507391030006
$ws.Range("O15").Value = 0.08457024278578675
$ws.Range("P15").Value = 0.08457024278578675
$ws.Range("Q15").Value = 6.295360846315999
$ws.Range("R15").Value = 56.65824761684399
$ws.Range("S15").Value = 0.00105502002936471
$ws.Range("T15").Value = 0.00105502002936471

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6381486666666666
$ws.Range("H16").Value = 1.914446
$ws.Range("I16").Value = 0.01247507391030006
$ws.Range("J16").Value = 0.01247507391030006
$ws.Range("M16").Value = 35.05835333333334
$ws.Range("N16").Value = 105.17506
$ws.Range("O16").Value = 0.3005455684073286
$ws.Range("P16").Value = 0.3005455684073286
$ws.Range("Q16").Value = 22.37244143519555
$ws.Range("R16").Value = 201.35197291676
$ws.Range("S16").Value = 0.003749328179294568
$ws.Range("T16").Value = 0.003749328179294568

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6381486666666666
$ws.Range("H17").Value = 1.914446
$ws.Range("I17").Value = 0.01247507391030006
$ws.Range("J17").Value = 0.01247507391030006
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.183961
$ws.Range("N17").Value = 0.551883
$ws.Range("O17").Value = 0.001577046782092083
$ws.Range("P17").Value = 0.001577046782092083
$ws.Range("Q17").Value = 0.1173944668686667
$ws.Range("R17").Value = 1.056550201818
$ws.Range("S17").Value = 0.00001967377516659961
$ws.Range("T17").Value = 0.00001967377516659961
